# CIERRE 17 SEPT 22
# Advance the payroll week from "SEMANA 36 (05-11 SEPT 2022)" to
# "SEMANA 37 (12-18 SEPT 2022)" and update the figures for the new period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")

# --- Week header (drives H9/B27/H27/B43 via formulas) ---------------------
$ws.Range("B9").Value = "SEMANA  37  DEL    12      Al   18   DE  SEPTIEMBRE          2022"

# --- Employee #1 block (days / salary) -------------------------------------
$ws.Range("D3").Value = 5
$ws.Range("D3").Font.Bold = $true
$ws.Range("D3").Font.Size = 14
$ws.Range("D3").HorizontalAlignment = -4108   # xlCenter
$ws.Rows(3).RowHeight = 18.75

$ws.Range("E3").Value = 1833

$ws.Range("J4").Font.Bold = $true
$ws.Range("J4").Font.Size = 14
$ws.Range("J4").HorizontalAlignment = -4108   # xlCenter

# --- Employee #2 block (extras / totals) -----------------------------------
$ws.Range("K21").Value = 1540

$ws.Range("D22").Font.Bold = $true
$ws.Range("D22").Font.Size = 14
$ws.Range("D22").HorizontalAlignment = -4108  # xlCenter
$ws.Rows(22).RowHeight = 19.5

# --- Employee #5 block (extras) ---------------------------------------------
$ws.Range("E40").Value = 1250

# --- Refresh cached formula results (dates, sums, week-label copies) -------
$excel.Calculate()

# --- Restore the view: scroll position + active cell -----------------------
$ws.Range("E41").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
